# Eric Rogers.docx - "A Cat, a Parrot, and a Bag of Seed" problem
# Add an extra insight bullet before the (empty, bookmark-only) trailing
# list paragraph, and demote that trailing paragraph out of the list.

$d = $word.ActiveDocument

# The paragraph we need to split is the very last paragraph in the
# document body: it is empty except for the _GoBack bookmark, and is
# still part of the numbered/bulleted list (pStyle ListParagraph, numPr).
$lastPara = $d.Paragraphs.Last
$splitPos = $lastPara.Range.Start

# Insert a new paragraph mark right before that paragraph's content.
# Because the insertion point sits before the bookmark, the new empty
# paragraph created by the split ends up positioned first, while the
# original paragraph (and its bookmark) keeps its place as the last
# paragraph in the document.
$breakRange = $d.Range($splitPos, $splitPos)
$breakRange.InsertParagraphBefore()

# The freshly created paragraph (now second-to-last) inherited the list
# formatting automatically -- fill it in with the new insight text.
$paraCount = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($paraCount - 1)
$newPara.Range.InsertBefore("Something that isn’t explained in the word problem is that as the man takes each one to the other side and leaves it, he must keep in mind that he can’t leave the wrong two together on the other side while he gets the third.")

# The trailing (bookmark) paragraph should no longer be numbered -- only
# its ListParagraph style remains.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.ListFormat.RemoveNumbers()
